# Tutorial bullet currently reads:
#   "Welcome to Dark Water. Throughout the game you will be given
#    objectives by your helmsman "
# and should read:
#   "Welcome to Dark Water. Your goal in the game is to navigate the
#    ocean and complete the objectives given to you by your crew."
#
# Replace the run text ". Throughout the game you will be given
# objectives by your helmsman " (which spans four separate <w:r> runs)
# with the new sentence, leaving the preceding "Water" text and the
# trailing closing curly-quote untouched.

$d = $word.ActiveDocument

$old = ". Throughout the game you will be given objectives by your helmsman "
$new = ". Your goal in the game is to navigate the ocean and complete the objectives given to you by your crew."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the tutorial sentence to replace."
}

# Sanity-check the edit landed where expected (the tutorial bullet is the
# second paragraph of the outline).
$expectedTail = "Water. Your goal in the game is to navigate the ocean and complete the objectives given to you by your crew.”"
$actual = $d.Paragraphs(2).Range.Text
if ($actual -notlike "*$expectedTail*") {
    throw "Unexpected paragraph text after replace: $actual"
}
